$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D2"  = "307.44"
    "E2"  = "0.97%"
    "D3"  = "37.09"
    "E3"  = "3.85%"
    "D4"  = "5.048"
    "E4"  = "-1.16%"
    "D5"  = "0.07911"
    "E5"  = "0.86%"
    "D6"  = "2.203"
    "E6"  = "-3.54%"
    "D7"  = "8.021"
    "E7"  = "-1.38%"
    "B8"  = "MXToken"
    "C8"  = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
    "D8"  = "0.9300"
    "E8"  = "0.26%"
    "B9"  = "LiechtensteinCryptoassetsExchange"
    "C9"  = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
    "D9"  = "0.09920"
    "E9"  = "-0.75%"
    "B10" = "WazirX"
    "C10" = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
    "D10" = "0.1884"
    "E10" = "3.55%"
    "B11" = "MandalaExchangeToken"
    "C11" = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
    "D11" = "0.08704"
    "E11" = "0.38%"
    "B12" = "BitrueCoin"
    "C12" = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
    "D12" = "0.03599"
    "E12" = "6.61%"
    "B13" = "BitMartToken"
    "C13" = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
    "D13" = "0.09964"
    "E13" = "0.44%"
    "B14" = "BitForexToken"
    "C14" = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
    "D14" = "0.001482"
    "E14" = "0.78%"
    "B15" = "TigerCash"
    "C15" = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
    "D15" = "0.005635"
    "E15" = "0.15%"
    "B16" = "LEO"
    "C16" = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
    "D16" = "3.466"
    "E16" = "-0.52%"
    "B17" = "GateToken"
    "C17" = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
    "D17" = "4.032"
    "E17" = "0.30%"
    "D18" = "2.339"
    "E18" = "11.59%"
    "D19" = "0.3436"
    "E19" = "0.06%"
    "D20" = "0.1335"
    "E20" = "1.02%"
    "D21" = "4.942"
    "E21" = "8.96%"
    "E22" = "-1.79%"
    "D23" = "0.04587"
    "E23" = "-1.63%"
    "D24" = "0.005254"
    "E24" = "17.11%"
    "E25" = "0.81%"
    "D26" = "0.0001398"
    "E26" = "7.76%"
    "D27" = "0.0002713"
    "E27" = "0.39%"
    "D39" = "0.01836"
    "E39" = "4.17%"
    "D40" = "0.04786"
    "E40" = "1.56%"
    "D41" = "0.007990"
    "E41" = "1.84%"
    "D42" = "0.1416"
    "E42" = "0.05%"
    "E43" = "-11.76%"
    "D44" = "0.002111"
    "E44" = "-7.60%"
    "E45" = "13.37%"
    "D46" = "0.00006311"
    "E46" = "4.21%"
    "E47" = "-0.25%"
    "D48" = "0.0005791"
    "E48" = "-0.16%"
    "D49" = "35.73"
    "E49" = "516.66%"
    "D50" = "0.002685"
    "E50" = "-0.28%"
    "D51" = "0.00002097"
    "E51" = "-0.25%"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}
